$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.875.37'
$ws.Range("E2").Value = '  -1.24%  '

$ws.Range("D3").Value = '1.640.42'
$ws.Range("E3").Value = '  -0.87%  '

$fmt = $ws.Range("D4").NumberFormat
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").NumberFormat = $fmt
$ws.Range("E4").Value = '  -0.11%  '

$fmt = $ws.Range("D5").NumberFormat
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.58'
$ws.Range("D5").NumberFormat = $fmt
$ws.Range("E5").Value = '  -0.38%  '

$fmt = $ws.Range("D6").NumberFormat
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5024'
$ws.Range("D6").NumberFormat = $fmt
$ws.Range("E6").Value = '  -2.20%  '

$fmt = $ws.Range("D7").NumberFormat
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("D7").NumberFormat = $fmt
$ws.Range("E7").Value = '  -0.18%  '

$fmt = $ws.Range("D8").NumberFormat
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2569'
$ws.Range("D8").NumberFormat = $fmt
$ws.Range("E8").Value = '  -1.25%  '

$fmt = $ws.Range("D9").NumberFormat
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06369'
$ws.Range("D9").NumberFormat = $fmt
$ws.Range("E9").Value = '  -1.44%  '

$fmt = $ws.Range("D10").NumberFormat
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.57'
$ws.Range("D10").NumberFormat = $fmt
$ws.Range("E10").Value = '  -1.88%  '

$fmt = $ws.Range("D11").NumberFormat
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07758'
$ws.Range("D11").NumberFormat = $fmt
$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").Value = '1.652.65'
$ws.Range("E12").Value = '  -0.16%  '

$fmt = $ws.Range("D13").NumberFormat
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.249'
$ws.Range("D13").NumberFormat = $fmt
$ws.Range("E13").Value = '  -1.28%  '

$ws.Range("D14").Value = '1.865.81'
$ws.Range("E14").Value = '  -0.92%  '

$fmt = $ws.Range("D15").NumberFormat
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5453'
$ws.Range("D15").NumberFormat = $fmt
$ws.Range("E15").Value = '  -1.43%  '

$ws.Range("D16").Value = '0.0₅7875'
$ws.Range("E16").Value = '  -2.12%  '

$fmt = $ws.Range("D17").NumberFormat
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.01'
$ws.Range("D17").NumberFormat = $fmt
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").Value = '25.900.02'
$ws.Range("E18").Value = '  -1.30%  '

$ws.Range("E19").Value = '  -0.14%  '

$fmt = $ws.Range("D20").NumberFormat
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.85'
$ws.Range("D20").NumberFormat = $fmt

$fmt = $ws.Range("D21").NumberFormat
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.381'
$ws.Range("D21").NumberFormat = $fmt
$ws.Range("E21").Value = '  -0.93%  '

$fmt = $ws.Range("D22").NumberFormat
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.876'
$ws.Range("D22").NumberFormat = $fmt
$ws.Range("E22").Value = '  -2.10%  '

$fmt = $ws.Range("D23").NumberFormat
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.967'
$ws.Range("D23").NumberFormat = $fmt
$ws.Range("E23").Value = '  -1.02%  '

$ws.Range("E24").Value = '  -0.13%  '

$fmt = $ws.Range("D25").NumberFormat
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.892'
$ws.Range("D25").NumberFormat = $fmt
$ws.Range("E25").Value = '  +5.02%  '

$ws.Range("E26").Value = '  -2.89%  '

$ws.Range("E27").Value = '  -4.02%  '

$fmt = $ws.Range("D28").NumberFormat
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.64'
$ws.Range("D28").NumberFormat = $fmt
$ws.Range("E28").Value = '  -1.43%  '

$fmt = $ws.Range("D29").NumberFormat
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.761'
$ws.Range("D29").NumberFormat = $fmt
$ws.Range("E29").Value = '  -3.68%  '

$ws.Range("E30").Value = '  -0.09%  '

$fmt = $ws.Range("D31").NumberFormat
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04956'
$ws.Range("D31").NumberFormat = $fmt
$ws.Range("E31").Value = '  -2.97%  '

$fmt = $ws.Range("D32").NumberFormat
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.253'
$ws.Range("D32").NumberFormat = $fmt
$ws.Range("E32").Value = '  -3.36%  '

$fmt = $ws.Range("D33").NumberFormat
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.189'
$ws.Range("D33").NumberFormat = $fmt
$ws.Range("E33").Value = '  -1.40%  '

$fmt = $ws.Range("D34").NumberFormat
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.543'
$ws.Range("D34").NumberFormat = $fmt
$ws.Range("E34").Value = '  -1.42%  '

$ws.Range("E35").Value = '  +0.82%  '

$fmt = $ws.Range("D36").NumberFormat
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.631'
$ws.Range("D36").NumberFormat = $fmt
$ws.Range("E36").Value = '  -3.76%  '

$fmt = $ws.Range("D37").NumberFormat
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8902'
$ws.Range("D37").NumberFormat = $fmt
$ws.Range("E37").Value = '  -3.85%  '

$fmt = $ws.Range("D38").NumberFormat
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5620'
$ws.Range("D38").NumberFormat = $fmt
$ws.Range("E38").Value = '  -2.07%  '

$ws.Range("D39").Value = '1.146.37'
$ws.Range("E39").Value = '  -1.67%  '

$fmt = $ws.Range("D40").NumberFormat
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01565'
$ws.Range("D40").NumberFormat = $fmt
$ws.Range("E40").Value = '  -1.61%  '

$ws.Range("E41").Value = '  -0.11%  '

$fmt = $ws.Range("D42").NumberFormat
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.664'
$ws.Range("D42").NumberFormat = $fmt
$ws.Range("E42").Value = '  -0.94%  '

$fmt = $ws.Range("D43").NumberFormat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.87'
$ws.Range("D43").NumberFormat = $fmt
$ws.Range("E43").Value = '  -0.51%  '

$fmt = $ws.Range("D44").NumberFormat
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8058'
$ws.Range("D44").NumberFormat = $fmt
$ws.Range("E44").Value = '  -2.28%  '

$ws.Range("D45").Value = '1.776.87'

$ws.Range("D46").Value = '0.0₈116'
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("E47").Value = '  +0.03%  '

$fmt = $ws.Range("D48").NumberFormat
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("D48").NumberFormat = $fmt
$ws.Range("E48").Value = '  -0.07%  '

$fmt = $ws.Range("D49").NumberFormat
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.74'
$ws.Range("D49").NumberFormat = $fmt
$ws.Range("E49").Value = '  -1.32%  '

$fmt = $ws.Range("D50").NumberFormat
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05057'
$ws.Range("D50").NumberFormat = $fmt

$ws.Range("E51").Value = '  -0.25%  '
